$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Protocol")

# Update the "you_close_it" comment (row 16, column F) with the new, longer
# description of the active-close behavior.
$ws.Range("F16").Value2 = "This is useful when one side wants the other to initiate the active close on the TCP socket for this transport (usually, to avoid having TIME_WAIT sockets). Only server can ask client to do this. If you_close_it is sent over an HTTP transport, client must try to close the HTTP connection. (usually with xhrObject.abort() or removing an iframe)"

# Update the "sack, seqNum, sackedList" comment (row 12, column F) to the
# new wording about freeing memory in the peer's box queue.
$ws.Range("F12").Value2 = "Both parties need to sack often to free memory in their peer's box queue"

# The longer you_close_it comment needs extra vertical room.
$ws.Rows.Item(16).RowHeight = 40.5

# Move the view: scroll back to the top (clear any frozen/scrolled top-left
# cell) and move the active selection to A12.
$ws.Range("A12").Select()

$wb.Save()
